$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new column ("culture_collection") is being inserted before column AH (34).
# Row 15 (the header row) is the only row with data beyond column A, so shift
# its cells manually instead of doing a full EntireColumn insert (which would also
# alter the pre-existing <cols> width definitions that otherwise stay untouched).

# First, clone formatting (style) from an existing data-header cell into the brand
# new rightmost cell (column 97) that will receive the shifted "water_content" header.
$ws.Cells.Item(15, 19).Copy($ws.Cells.Item(15, 97))

# Shift header text one column to the right, from the rightmost column down to AH,
# so no values get overwritten before they are copied onward.
for ($col = 96; $col -ge 34; $col--) {
    $srcCell = $ws.Cells.Item(15, $col)
    $dstCell = $ws.Cells.Item(15, $col + 1)
    $dstCell.Value = $srcCell.Value()
}

# Comments stay attached to their original cell during the manual value shift above,
# so move each comment one column to the right as well (rightmost first).
$srcCell = $ws.Cells.Item(15, 96)
$dstCell = $ws.Cells.Item(15, 97)
$dstCell.AddComment('water content measurement')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 95)
$dstCell = $ws.Cells.Item(15, 96)
$dstCell.AddComment('turbidity measurement')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 94)
$dstCell = $ws.Cells.Item(15, 95)
$dstCell.AddComment('Feeding position in food chain (eg., chemolithotroph)')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 93)
$dstCell = $ws.Cells.Item(15, 94)
$dstCell.AddComment('Definition for soil: total organic C content of the soil units of g C/kg soil. Definition otherwise: total organic carbon content')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 92)
$dstCell = $ws.Cells.Item(15, 93)
$dstCell.AddComment('total nitrogen content of the sample')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 91)
$dstCell = $ws.Cells.Item(15, 92)
$dstCell.AddComment('total carbon content')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 90)
$dstCell = $ws.Cells.Item(15, 91)
$dstCell.AddComment('temperature of the sample at time of sampling')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 89)
$dstCell = $ws.Cells.Item(15, 90)
$dstCell.AddComment('concentration of sulfide')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 88)
$dstCell = $ws.Cells.Item(15, 89)
$dstCell.AddComment('concentration of sulfate')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 87)
$dstCell = $ws.Cells.Item(15, 88)
$dstCell.AddComment('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 86)
$dstCell = $ws.Cells.Item(15, 87)
$dstCell.AddComment('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 85)
$dstCell = $ws.Cells.Item(15, 86)
$dstCell.AddComment('sodium concentration')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 84)
$dstCell = $ws.Cells.Item(15, 85)
$dstCell.AddComment('concentration of silicate')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 83)
$dstCell = $ws.Cells.Item(15, 84)
$dstCell.AddComment('volume (mL) or weight (g) of sample processed for DNA extraction')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 82)
$dstCell = $ws.Cells.Item(15, 83)
$dstCell.AddComment('temperature at which sample was stored, e.g. -80')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 81)
$dstCell = $ws.Cells.Item(15, 82)
$dstCell.AddComment('location at which sample was stored, usually name of a specific freezer/room')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 80)
$dstCell = $ws.Cells.Item(15, 81)
$dstCell.AddComment('duration for which sample was stored')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 79)
$dstCell = $ws.Cells.Item(15, 80)
$dstCell.AddComment('Amount or size of sample (volume, mass or area) that was collected')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 78)
$dstCell = $ws.Cells.Item(15, 79)
$dstCell.AddComment('Processing applied to the sample during or after isolation')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 77)
$dstCell = $ws.Cells.Item(15, 78)
$dstCell.AddComment('Method or device employed for collecting sample')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 76)
$dstCell = $ws.Cells.Item(15, 77)
$dstCell.AddComment('salinity measurement')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 75)
$dstCell = $ws.Cells.Item(15, 76)
$dstCell.AddComment('Aerobic or anaerobic')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 74)
$dstCell = $ws.Cells.Item(15, 75)
$dstCell.AddComment('redox potential, measured relative to a hydrogen cell, indicating oxidation or reduction potential')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 73)
$dstCell = $ws.Cells.Item(15, 74)
$dstCell.AddComment('pressure to which the sample is subject, in atmospheres')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 72)
$dstCell = $ws.Cells.Item(15, 73)
$dstCell.AddComment('concentration of potassium')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 71)
$dstCell = $ws.Cells.Item(15, 72)
$dstCell.AddComment('concentration of phospholipid fatty acids; can include multiple values')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 70)
$dstCell = $ws.Cells.Item(15, 71)
$dstCell.AddComment('concentration of phosphate')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 69)
$dstCell = $ws.Cells.Item(15, 70)
$dstCell.AddComment('concentration of phaeopigments; can include multiple phaeopigments')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 68)
$dstCell = $ws.Cells.Item(15, 69)
$dstCell.AddComment('pH measurement')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 67)
$dstCell = $ws.Cells.Item(15, 68)
$dstCell.AddComment('concentration of petroleum hydrocarbon')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 66)
$dstCell = $ws.Cells.Item(15, 67)
$dstCell.AddComment('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 65)
$dstCell = $ws.Cells.Item(15, 66)
$dstCell.AddComment('To what is the entity pathogenic')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 64)
$dstCell = $ws.Cells.Item(15, 65)
$dstCell.AddComment('concentration of particulate organic carbon')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 63)
$dstCell = $ws.Cells.Item(15, 64)
$dstCell.AddComment('oxygenation status of sample')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 62)
$dstCell = $ws.Cells.Item(15, 63)
$dstCell.AddComment('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 61)
$dstCell = $ws.Cells.Item(15, 62)
$dstCell.AddComment('concentration of organic nitrogen')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 60)
$dstCell = $ws.Cells.Item(15, 61)
$dstCell.AddComment('concentration of organic matter')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 59)
$dstCell = $ws.Cells.Item(15, 60)
$dstCell.AddComment('concentration of organic carbon')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 58)
$dstCell = $ws.Cells.Item(15, 59)
$dstCell.AddComment('concentration of nitrogen (total)')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 57)
$dstCell = $ws.Cells.Item(15, 58)
$dstCell.AddComment('concentration of nitrite')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 56)
$dstCell = $ws.Cells.Item(15, 57)
$dstCell.AddComment('concentration of nitrate')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 55)
$dstCell = $ws.Cells.Item(15, 56)
$dstCell.AddComment('concentration of n-alkanes; can include multiple n-alkanes')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 54)
$dstCell = $ws.Cells.Item(15, 55)
$dstCell.AddComment('any other measurement performed or parameter collected, that is not listed here')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 53)
$dstCell = $ws.Cells.Item(15, 54)
$dstCell.AddComment('methane (gas) amount or concentration at the time of sampling')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 52)
$dstCell = $ws.Cells.Item(15, 53)
$dstCell.AddComment('measurement of mean peak friction velocity')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 51)
$dstCell = $ws.Cells.Item(15, 52)
$dstCell.AddComment('measurement of mean friction velocity')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 50)
$dstCell = $ws.Cells.Item(15, 51)
$dstCell.AddComment('concentration of magnesium')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 49)
$dstCell = $ws.Cells.Item(15, 50)
$dstCell.AddComment('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 48)
$dstCell = $ws.Cells.Item(15, 49)
$dstCell.AddComment('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 47)
$dstCell = $ws.Cells.Item(15, 48)
$dstCell.AddComment('NCBI taxonomy ID of the host, e.g. 9606')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 46)
$dstCell = $ws.Cells.Item(15, 47)
$dstCell.AddComment('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 45)
$dstCell = $ws.Cells.Item(15, 46)
$dstCell.AddComment('Health or disease status of sample at time of collection')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 44)
$dstCell = $ws.Cells.Item(15, 45)
$dstCell.AddComment('measurement of glucosidase activity')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 43)
$dstCell = $ws.Cells.Item(15, 44)
$dstCell.AddComment('Plasmids that have significance phenotypic consequence')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 42)
$dstCell = $ws.Cells.Item(15, 43)
$dstCell.AddComment('Estimated size of genome')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 41)
$dstCell = $ws.Cells.Item(15, 42)
$dstCell.AddComment('Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 40)
$dstCell = $ws.Cells.Item(15, 41)
$dstCell.AddComment('concentration of dissolved oxygen')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 39)
$dstCell = $ws.Cells.Item(15, 40)
$dstCell.AddComment('dissolved organic nitrogen concentration measured as; total dissolved nitrogen - NH4 - NO3 - NO2')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 38)
$dstCell = $ws.Cells.Item(15, 39)
$dstCell.AddComment('concentration of dissolved organic carbon')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 37)
$dstCell = $ws.Cells.Item(15, 38)
$dstCell.AddComment('dissolved inorganic carbon concentration')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 36)
$dstCell = $ws.Cells.Item(15, 37)
$dstCell.AddComment('concentration of dissolved hydrogen')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 35)
$dstCell = $ws.Cells.Item(15, 36)
$dstCell.AddComment('concentration of dissolved carbon dioxide')
$srcCell.Comment.Delete()

$srcCell = $ws.Cells.Item(15, 34)
$dstCell = $ws.Cells.Item(15, 35)
$dstCell.AddComment('concentration of diether lipids; can include multiple types of diether lipids')
$srcCell.Comment.Delete()

# Finally, set the new "culture_collection" header text and its comment at AH15 (column 34).
$ws.Cells.Item(15, 34).Value = "culture_collection"
$ws.Cells.Item(15, 34).AddComment('Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier')
